$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Title text: " AGENDA Barroc-IT Groep" -> " AGENDA Barroc-IT Groep 5"
$ws.Range("E3").Value = " AGENDA Barroc-IT Groep 5"

# Date text: "Datum:25 aug." -> "Datum: 8 sep."
$ws.Range("B5").Value = "Datum: 8 sep."

# Minutes durations added for the first few agenda items
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 10
$ws.Range("C16").Value = 5

# D10: "k" -> "d"
$ws.Range("D10").Value = "d"

# F10: new text
$ws.Range("F10").Value = "Planning (vooruitzicht)"

# Reflect the author's final view state (scroll position, zoom, selection)
$ws.Activate()
[void]$ws.Range("F25").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 168
